$d = $word.ActiveDocument

# Remove the trailing "Requisitos" heading paragraph and the following
# "LOB1045 - ..." (ListBullet) paragraph that lists the course requisite.
$count = $d.Paragraphs.Count
$startPara = $d.Paragraphs.Item($count - 1)
$endPara = $d.Paragraphs.Item($count)

$rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
$rng.Delete()
